# Apply the commit's changes:
#  - update c2db_id values for the three "None" placeholder rows (B2 stays same index,
#    B4 and B5 get distinct new ids)
#  - clear the leftover empty "neighbour" marker column (G) for rows 64-122
#  - bump the header/data row heights (rows 1-47) from 17.25 to 18.75
#  - normalize the font color used by the header/text style from theme color to explicit black
#  - drop the border/quote-prefix formatting that used to back column G's marker style,
#    matching the new (unstyled) look of columns A/B

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the c2db_id text values
$ws.Range("B2").Value = "e50c9fb190b0"
$ws.Range("B4").Value = "276f0a298324"
$ws.Range("B5").Value = "3deec82af6d4"

# 2) Clear the stray empty marker cells in column G (rows 64 through 122) and make
#    them match the plain (borderless) look of columns A/B instead of the old
#    bordered/quote-prefixed marker style
$gRange = $ws.Range("G64:G122")
$gRange.ClearContents()
$gRange.Borders.LineStyle = 0
$gRange.HorizontalAlignment = -4131   # xlLeft
$gRange.Font.Name = "Calibri"
$gRange.Font.Size = 11

# 3) Increase row heights for rows 1-47
for ($r = 1; $r -le 47; $r++) {
    $ws.Rows.Item($r).RowHeight = 18.75
}

# 4) Set explicit black font color (was theme color) on the styled text range
#    (columns A-C, which carry the font whose color definition changes; this also
#    covers the now-plain column G cells once they share that same style)
$textRange = $ws.Range("A1:C125")
$textRange.Font.Color = 0
$gRange.Font.Color = 0

